$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71, column A: timestamp value was re-recorded with slightly higher precision.
$ws.Range("A71").Value = 44384.76911854051

# New row 72: a freshly retrieved data row appended to the log.
$ws.Range("A72").Value = 44385.77016085689
$ws.Range("B72").Value = 80018
$ws.Range("C72").Value = 67494
$ws.Range("D72").Value = 3587
$ws.Range("E72").Value = 2196
$ws.Range("F72").Value = 1568
$ws.Range("G72").Value = 21233
$ws.Range("H72").Value = 1565
$ws.Range("I72").Value = 885
$ws.Range("J72").Value = 199
